{"js": "// Add a new \"Testing on 30th April 2023\" BodyText paragraph right after\n// the existing \"Testing on 16th March 2023\" paragraph (just before the\n// \"hello\" bookmark's end), mirroring the sibling date-testing paragraphs\n// already in the document.\nconst body = context.document.body;\n\nconst results = body.search(\"Testing on 16th March 2023\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find paragraph \"Testing on 16th March 2023\"');\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"Testing on 30th April 2023\",\n  \"After\"\n);\nnewParagraph.style = \"Body Text\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Confirm the anchor paragraph exists (mirrors how this edit would be\n# located interactively with Find) before touching the document.\n$found = $d.Content.Find.Execute(\"Testing on 16th March 2023\")\nif (-not $found) {\n    throw 'Could not find paragraph \"Testing on 16th March 2023\"'\n}\n\n# Re-resolve the anchor paragraph through the Paragraphs collection by\n# matching text: Range handles returned by Find can go stale once the\n# document is structurally edited, so look it up by index instead of\n# holding on to the Find range across the InsertParagraphAfter call.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n    if ($paraText -eq \"Testing on 16th March 2023\") {\n        $anchorIndex = $i\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not resolve anchor paragraph index\"\n}\n\n$anchorParagraph = $d.Paragraphs.Item($anchorIndex)\n$anchorParagraph.Range.InsertParagraphAfter()\n\n# The freshly inserted paragraph is now immediately after the anchor;\n# give it the same BodyText style and the new sentence.\n$newParagraph = $d.Paragraphs.Item($anchorIndex + 1)\n$newParagraph.Range.Text = \"Testing on 30th April 2023\"\n$newParagraph.Style = \"Body Text\"\n"}
